# GroupByTest_desc.xlsx edit
# - jx:each on A2: drop the var="g" attribute (groupBy no longer needs a var prefix)
# - jx:each on A3: items="g.items" -> items="_group.items"
# - A2 cell text: ${g.item.salaryGroup} -> ${_group.item.salaryGroup}
# - Active selection moves from A4 to A13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Comment on A2 (jx:each for the groupBy loop) - remove var="g"
$ws.Range("A2").Comment.Text("jxlsteam:`njx:each(items=""employees"" groupBy=""salaryGroup"" groupOrder=""DESC"" lastCell=""C3"")") | Out-Null

# Comment on A3 (jx:each for the items inside a group) - items="g.items" -> items="_group.items"
$ws.Range("A3").Comment.Text("jxlsteam:`njx:each(items=""_group.items"" var=""e"" lastCell=""C3"")") | Out-Null

# Cell A2 text uses the implicit group variable now
$ws.Range("A2").Value = '${_group.item.salaryGroup}'

# Selection moved to A13
$ws.Range("A13").Select() | Out-Null
